$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 457.33334
$ws.Range("I6").Value = 190
$ws.Range("J6").Value = 992
$ws.Range("K6").Value = 570
$ws.Range("L6").Value = 2976
$ws.Range("M6").Value = -458
$ws.Range("N6").Value = -3200
$ws.Range("H9").Value = 139.66667
$ws.Range("I9").Value = 139.66667
$ws.Range("K9").Value = 139.66667
$ws.Range("M9").Value = 29.33332999999999
$ws.Range("H12").Value = 111113016
$ws.Range("I12").Value = 2028.5
$ws.Range("J12").Value = 333335000
$ws.Range("K12").Value = 2028.5
$ws.Range("L12").Value = 333335000
$ws.Range("M12").Value = -1858.5
$ws.Range("N12").Value = -333335340
$ws.Range("H21").Value = 7839.4443
$ws.Range("I21").Value = 6106.8
$ws.Range("J21").Value = 10005.25
$ws.Range("K21").Value = 6106.8
$ws.Range("L21").Value = 10005.25
$ws.Range("M21").Value = -5638.8
$ws.Range("N21").Value = -10941.25
$ws.Range("H23").Value = 7839.4443
$ws.Range("I23").Value = 6106.8
$ws.Range("J23").Value = 10005.25
$ws.Range("K23").Value = 6106.8
$ws.Range("L23").Value = 10005.25
$ws.Range("M23").Value = -5872.8
$ws.Range("N23").Value = -10473.25
$ws.Range("H29").Value = 1061
$ws.Range("J29").Value = 4004
$ws.Range("L29").Value = 12012
$ws.Range("N29").Value = -12574
$ws.Range("H38").Value = 320.90475
$ws.Range("I38").Value = 136.95
$ws.Range("J38").Value = 4000
$ws.Range("K38").Value = 410.85
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = -38.84999999999997
$ws.Range("N38").Value = -12744
$ws.Range("H58").Value = 1458.8572
$ws.Range("J58").Value = 4017
$ws.Range("L58").Value = 12051
$ws.Range("N58").Value = -12351
$ws.Range("H80").Value = 5048.5557
$ws.Range("I80").Value = 2434.7273
$ws.Range("J80").Value = 6845.5625
$ws.Range("K80").Value = 7304.1819
$ws.Range("L80").Value = 20536.6875
$ws.Range("M80").Value = -6306.1819
$ws.Range("N80").Value = -22532.6875
$ws.Range("H83").Value = 5048.5557
$ws.Range("I83").Value = 2434.7273
$ws.Range("J83").Value = 6845.5625
$ws.Range("K83").Value = 21912.5457
$ws.Range("L83").Value = 61610.0625
$ws.Range("M83").Value = -16920.5457
$ws.Range("N83").Value = -71594.0625
$ws.Range("H137").Value = 2022.3462
$ws.Range("I137").Value = 1794.6316
$ws.Range("J137").Value = 2640.4285
$ws.Range("K137").Value = 5383.8948
$ws.Range("L137").Value = 7921.2855
$ws.Range("M137").Value = -2833.8948
$ws.Range("N137").Value = -13021.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 43679.848
$ws.Range("I74").Value = 51317.5
$ws.Range("J74").Value = 18221
$ws.Range("K74").Value = 51317.5
$ws.Range("L74").Value = 18221
$ws.Range("M74").Value = -50443.5
$ws.Range("N74").Value = -19969
$ws.Range("H77").Value = 43679.848
$ws.Range("I77").Value = 51317.5
$ws.Range("J77").Value = 18221
$ws.Range("K77").Value = 256587.5
$ws.Range("L77").Value = 91105
$ws.Range("M77").Value = -252219.5
$ws.Range("N77").Value = -99841

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2285.3845
$ws.Range("I31").Value = 2074.2856
$ws.Range("J31").Value = 2531.6667
$ws.Range("K31").Value = 2074.2856
$ws.Range("L31").Value = 2531.6667
$ws.Range("M31").Value = -1779.2856
$ws.Range("N31").Value = -3121.6667
$ws.Range("H34").Value = 2285.3845
$ws.Range("I34").Value = 2074.2856
$ws.Range("J34").Value = 2531.6667
$ws.Range("K34").Value = 2074.2856
$ws.Range("L34").Value = 2531.6667
$ws.Range("M34").Value = -1872.2856
$ws.Range("N34").Value = -2935.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 390.57144
$ws.Range("I5").Value = 329.44
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 988.3199999999999
$ws.Range("L5").Value = 2700
$ws.Range("M5").Value = -876.3199999999999
$ws.Range("N5").Value = -2924
$ws.Range("H17").Value = 597.1667
$ws.Range("I17").Value = 172.42857
$ws.Range("J17").Value = 1191.8
$ws.Range("K17").Value = 517.28571
$ws.Range("L17").Value = 3575.4
$ws.Range("M17").Value = -348.28571
$ws.Range("N17").Value = -3913.4
$ws.Range("H34").Value = 300.75
$ws.Range("I34").Value = 300.75
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 902.25
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -818.25
$ws.Range("N34").ClearContents()
$ws.Range("H39").Value = 1921.4286
$ws.Range("I39").Value = 450
$ws.Range("J39").Value = 3025
$ws.Range("K39").Value = 1350
$ws.Range("L39").Value = 9075
$ws.Range("M39").Value = -1056
$ws.Range("N39").Value = -9663
$ws.Range("H55").Value = 976.5
$ws.Range("I55").Value = 377
$ws.Range("J55").Value = 1276.25
$ws.Range("K55").Value = 1131
$ws.Range("L55").Value = 3828.75
$ws.Range("M55").Value = -954
$ws.Range("N55").Value = -4182.75
$ws.Range("H122").Value = 1224.9474
$ws.Range("I122").Value = 549.75
$ws.Range("J122").Value = 1716
$ws.Range("K122").Value = 4947.75
$ws.Range("L122").Value = 15444
$ws.Range("M122").Value = -2497.75
$ws.Range("N122").Value = -20344
$ws.Range("H132").Value = 3126.2173
$ws.Range("I132").Value = 2383.2222
$ws.Range("J132").Value = 5801
$ws.Range("K132").Value = 21448.9998
$ws.Range("L132").Value = 52209
$ws.Range("M132").Value = -18918.9998
$ws.Range("N132").Value = -57269
$ws.Range("H135").Value = 390.57144
$ws.Range("I135").Value = 329.44
$ws.Range("J135").Value = 900
$ws.Range("K135").Value = 2964.96
$ws.Range("L135").Value = 8100
$ws.Range("M135").Value = -429.96
$ws.Range("N135").Value = -13170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 17571.857
$ws.Range("I68").Value = 35666.668
$ws.Range("J68").Value = 4000.75
$ws.Range("K68").Value = 35666.668
$ws.Range("L68").Value = 4000.75
$ws.Range("M68").Value = -34917.668
$ws.Range("N68").Value = -5498.75
$ws.Range("H71").Value = 17571.857
$ws.Range("I71").Value = 35666.668
$ws.Range("J71").Value = 4000.75
$ws.Range("K71").Value = 178333.34
$ws.Range("L71").Value = 20003.75
$ws.Range("M71").Value = -174589.34
$ws.Range("N71").Value = -27491.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 33593.332
$ws.Range("J128").Value = 33593.332
$ws.Range("L128").Value = 33593.332
$ws.Range("N128").Value = -43553.332

Write-Host "Done applying edits"
